$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held 4 waves of Nigeria (GHS) data in rows 2-5.
# This update replaces it with a single row of Cambodia (CAS) data,
# dropping the now-obsolete wave rows 3-5 entirely.
$ws.Rows("3:5").Delete()

# Overwrite the remaining data row (row 2) with the Cambodia figures.
$ws.Range("A2").Value = "CAS"
$ws.Range("B2").Value = "khm"
$ws.Range("C2").Value = "🇰🇭"
$ws.Range("D2").Value = "w1"
$ws.Range("E2").Value = "Cambodia"
$ws.Range("F2").Value = 2019
$ws.Range("G2").Value = "2019-2020"
# H2:L2 (inflation, gdp_ppp, cons_ppp, gdp_ppp_2017, cons_ppp_2017) are left
# unchanged - the underlying figures carry over as-is.

# Match the author's final selection.
$ws.Range("G2").Select()
